# Credit.xlsx — "blender tutorial & SHanpe moving" commit
# Sheet1 ("Referensi"): insert a new "Pro Builder" reference row right after the
# header, and append a new "Blender Fundamentals" reference row at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new row 2 (old rows 2-7 shift down to 3-8) ------------------
$ws.Rows.Item(2).Insert()

# New row 2: Pro Builder intro video (Brackeys).
# Written in this column order (B, C, E, D) so new shared-strings are interned
# in the same sequence as the authored workbook.
$ws.Range("A2").Value = "Brackeys"
$ws.Range("B2").Value = "https://www.youtube.com/watch?v=PUSOg5YEflM"
$ws.Range("C2").Value = "Pro Builder"
$ws.Range("E2").Value = " "
$ws.Range("D2").Value = "Pro Builder Perkenalan tutorial"

# New row 9 (appended after the previously-last row, now row 8): Blender Fundamentals
$ws.Range("A9").Value = "Blender Foundation"
$ws.Range("B9").Value = "https://www.youtube.com/playlist?list=PLa1F2ddGya_8V90Kd5eC5PeBjySbXWGK1"
$ws.Range("C9").Value = "Blender Fundamentals"
$ws.Range("D9").Value = "Blender tutorial beginner lengkap"

# --- Hyperlinks -------------------------------------------------------------
# The row insert shifted the old NavMesh-video hyperlink's data down to B3, but
# the link object itself doesn't follow automatically, so rebuild both links.
$ws.Range("A1:E9").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B3"), "https://www.youtube.com/watch?v=CHV1ymlw-P8") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.youtube.com/watch?v=PUSOg5YEflM", [Type]::Missing, [Type]::Missing, "https://www.bing.com/search?q=unity+probuilder&form=EDGTCT&qs=AS&cvid=21509d34d4484303becc906eef91a8fb&refig=a8e419fc3ce1430fe0c1388e8b2d3902&cc=ID&setlang=id-ID&plvar=0&PC=ASTS") | Out-Null
# Adding the hyperlink with display text overwrites the cell's own text with
# that display text, so restore the real URL as the cell's stored value.
$ws.Range("B2").Value = "https://www.youtube.com/watch?v=PUSOg5YEflM"

# Restore the standard "Hipertaut" hyperlink look on both linked cells.
$ws.Range("B2").Style = "Hipertaut"
$ws.Range("B3").Style = "Hipertaut"

# --- View / selection --------------------------------------------------------
# Moves the active cell (and drops the old topLeftCell scroll anchor).
$ws.Range("D9").Select() | Out-Null

# --- Page setup ---------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
